$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "64.130.23"
$ws.Cells.Item(2, 5).Value = "  +2.24%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "2.510.89"
$ws.Cells.Item(3, 5).Value = "  +2.13%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "580.87"
$ws.Cells.Item(5, 5).Value = "  +1.28%  "

# Row 6 - Solana
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "151.90"
$ws.Cells.Item(6, 5).Value = "  +4.52%  "

# Row 7 - USDC
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.00"
$ws.Cells.Item(7, 5).Value = "  -0.09%  "

# Row 8 - XRP
$ws.Cells.Item(8, 5).Value = "  +0.39%  "

# Row 9 - Dogecoin
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.113"
$ws.Cells.Item(9, 5).Value = "  +1.51%  "

# Row 10 - TRON
$ws.Cells.Item(10, 5).Value = "  +0.22%  "

# Row 11 - Toncoin
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.26"
$ws.Cells.Item(11, 5).Value = "  +0.47%  "

# Row 12 - Cardano
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.355"
$ws.Cells.Item(12, 5).Value = "  -0.12%  "

# Row 13 - Avalanche
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "29.73"
$ws.Cells.Item(13, 5).Value = "  +3.19%  "

# Row 14 - ShibaInu
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.0000180"
$ws.Cells.Item(14, 5).Value = "  +2.03%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Cells.Item(15, 4).Value = "2.964.66"
$ws.Cells.Item(15, 5).Value = "  +2.04%  "

# Row 16 - WrappedBTC
$ws.Cells.Item(16, 4).Value = "63.182.37"
$ws.Cells.Item(16, 5).Value = "  +0.77%  "

# Row 17 - WrappedEther
$ws.Cells.Item(17, 4).Value = "2.509.07"
$ws.Cells.Item(17, 5).Value = "  +1.82%  "

# Row 18 - Uniswap
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "7.87"
$ws.Cells.Item(18, 5).Value = "  -1.50%  "

# Row 19 - Chainlink
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "11.00"
$ws.Cells.Item(19, 5).Value = "  +0.39%  "

# Row 20 - Polkadot
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "4.25"
$ws.Cells.Item(20, 5).Value = "  +2.82%  "

# Row 21 - SuiNetwork
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "2.27"
$ws.Cells.Item(21, 5).Value = "  +2.21%  "

# Row 22 - BitcoinCash
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "327.68"
$ws.Cells.Item(22, 5).Value = "  +0.44%  "

# Row 23 - Dai
$ws.Cells.Item(23, 5).Value = "  -0.01%  "

# Row 24 - Aptos
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "10.18"
$ws.Cells.Item(24, 5).Value = "  +1.41%  "

# Row 25 - Bittensor
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "673.01"
$ws.Cells.Item(25, 5).Value = "  +3.31%  "

# Row 26 - Litecoin
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "65.48"
$ws.Cells.Item(26, 5).Value = "  -0.13%  "

# Row 27 - PEPE
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.0000101"
$ws.Cells.Item(27, 5).Value = "  +3.24%  "

# Row 28 - WrappedeETH
$ws.Cells.Item(28, 5).Value = "  +0.65%  "

# Row 29 - Fetch.AI
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.50"
$ws.Cells.Item(29, 5).Value = "  +3.21%  "

# Row 30 - Binance-PegBSC-USD
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.988"
$ws.Cells.Item(30, 5).Value = "  -1.08%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "8.03"
$ws.Cells.Item(31, 5).Value = "  +0.72%  "

# Row 32 - PancakeSwap
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.86"
$ws.Cells.Item(32, 5).Value = "  +0.63%  "

# Row 33 - Kaspa
$ws.Cells.Item(33, 5).Value = "  +1.41%  "

# Row 34 - FirstDigitalUSD
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.998"
$ws.Cells.Item(34, 5).Value = "  -0.08%  "

# Row 35 - ImmutableX
$ws.Cells.Item(35, 5).Value = "  +1.82%  "

# Row 36 - NEARProtocol
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "4.83"
$ws.Cells.Item(36, 5).Value = "  +1.68%  "

# Row 37 - RenderToken
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "5.57"
$ws.Cells.Item(37, 5).Value = "  +4.04%  "

# Row 38 - PolygonEcosystemToken
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.371"
$ws.Cells.Item(38, 5).Value = "  +0.66%  "

# Row 39 - EthereumClassic (was Monero)
$ws.Cells.Item(39, 2).Value = "EthereumClassic"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "18.81"
$ws.Cells.Item(39, 5).Value = "  +0.83%  "

# Row 40 - Monero (was EthereumClassic)
$ws.Cells.Item(40, 2).Value = "Monero"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "152.06"
$ws.Cells.Item(40, 5).Value = "  -0.31%  "

# Row 41 - dogwifhat
$ws.Cells.Item(41, 5).Value = "  +1.91%  "

# Row 42 - Stacks
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.78"
$ws.Cells.Item(42, 5).Value = "  +3.66%  "

# Row 43 - USDe
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.999"
$ws.Cells.Item(43, 5).Value = "  +0.03%  "

# Row 44 - Aave
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "159.57"
$ws.Cells.Item(44, 5).Value = "  +4.26%  "

# Row 45 - BabyDogeCoin
$ws.Cells.Item(45, 4).Value = "0.0₆0301"
$ws.Cells.Item(45, 5).Value = "  -4.06%  "

# Row 46 - WhiteBITCoin
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "15.44"
$ws.Cells.Item(46, 5).Value = "  +1.45%  "

# Row 47 - Filecoin
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "3.63"
$ws.Cells.Item(47, 5).Value = "  +1.59%  "

# Row 48 - InjectiveProtocol
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "21.14"
$ws.Cells.Item(48, 5).Value = "  +4.51%  "

# Row 49 - Mantle
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.617"
$ws.Cells.Item(49, 5).Value = "  +1.95%  "

# Row 50 - Hedera
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0522"
$ws.Cells.Item(50, 5).Value = "  +2.23%  "

# Row 51 - VeChain (was Stellar)
$ws.Cells.Item(51, 2).Value = "VeChain"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0229"
$ws.Cells.Item(51, 5).Value = "  +1.81%  "
